$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New section header
$ws.Range("A11").Value = "With clamping FC layer:"
$ws.Range("A11").Font.Bold = $true

# Table header row (13)
$ws.Range("B13").Value = "fp32"
$ws.Range("C13").Value = "fp32"
$ws.Range("D13").Value = "ai84 quant"
$ws.Range("E13").Value = "ai85 quant"

# Sub-header row (14)
$ws.Range("A14").Value = "Dataset"
$ws.Range("A14").Font.Italic = $true
$ws.Range("B14").Value = "best (verif)"
$ws.Range("C14").Value = "test final"
$ws.Range("D14").Value = "test final"
$ws.Range("E14").Value = "test final"

# Data rows
$ws.Range("A15").Value = "MNIST"
$ws.Range("B15").Value = 99.4
$ws.Range("C15").Value = 99.6
$ws.Range("D15").Value = 99.5

$ws.Range("A16").Value = "FashionMNIST"
$ws.Range("B16").Value = 92.3
$ws.Range("C16").Value = 92.1
$ws.Range("D16").Value = 91.7

$ws.Range("A17").Value = "CIFAR-10"
$ws.Range("B17").Value = 82.6
$ws.Range("C17").Value = 82
$ws.Range("D17").Value = 82.4

$ws.Range("A18").Value = "CIFAR-10 w/bias"
$ws.Range("B18").Value = 82.7
$ws.Range("C18").Value = 82.1
$ws.Range("D18").Value = 31.1
$ws.Range("E18").Value = 81.599999999999994

$ws.Range("B15:D18").NumberFormat = "0.0"
$ws.Range("E18").NumberFormat = "0.0"

$ws.Range("E12").Select()
